$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "M1"
$ws.Range("D2").Value = "DO"
$ws.Range("I2").Value = "A1"
$ws.Range("J2").Value = "DO"
$ws.Range("N2").Value = "M3"
$ws.Range("P2").Value = "A1"
$ws.Range("Q2").Value = "A1"
$ws.Range("T2").Value = "DO"
$ws.Range("X2").Value = "DO"
$ws.Range("Y2").Value = "A1"
$ws.Range("Z2").Value = "A1"
$ws.Range("AB2").Value = "M1"
$ws.Range("AC2").Value = "M3"
$ws.Range("B3").Value = "A1"
$ws.Range("C3").Value = "A2"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = "M2"
$ws.Range("F3").Value = "DO"
$ws.Range("G3").Value = "M2"
$ws.Range("I3").Value = "DO"
$ws.Range("J3").Value = "M2"
$ws.Range("K3").Value = "A2"
$ws.Range("Q3").Value = "DO"
$ws.Range("S3").Value = "M2"
$ws.Range("V3").Value = "A1"
$ws.Range("W3").Value = "M2"
$ws.Range("X3").Value = "A2"
$ws.Range("Y3").Value = "A1"
$ws.Range("Z3").Value = "M2"
$ws.Range("AC3").Value = "DO"
$ws.Range("B4").Value = "DO"
$ws.Range("C4").Value = "M1"
$ws.Range("D4").Value = "A1"
$ws.Range("H4").Value = "M3"
$ws.Range("K4").Value = "A1"
$ws.Range("L4").Value = "A1"
$ws.Range("O4").Value = "A1"
$ws.Range("P4").Value = "DO"
$ws.Range("Q4").Value = "M1"
$ws.Range("S4").Value = "A1"
$ws.Range("T4").Value = "M1"
$ws.Range("U4").Value = "A1"
$ws.Range("V4").Value = "M3"
$ws.Range("AC4").Value = "A1"
$ws.Range("B5").Value = "M2"
$ws.Range("C5").Value = "M2"
$ws.Range("D5").Value = "DO"
$ws.Range("E5").Value = "M1"
$ws.Range("F5").Value = "A2"
$ws.Range("G5").Value = "M2"
$ws.Range("H5").Value = "A1"
$ws.Range("I5").Value = "M2"
$ws.Range("J5").Value = "A1"
$ws.Range("K5").Value = "A2"
$ws.Range("L5").Value = "DO"
$ws.Range("M5").Value = "M1"
$ws.Range("O5").Value = "M2"
$ws.Range("P5").Value = "M2"
$ws.Range("Q5").Value = "M2"
$ws.Range("R5").Value = "M2"
$ws.Range("S5").Value = "DO"
$ws.Range("U5").Value = "M1"
$ws.Range("X5").Value = "DO"
$ws.Range("Y5").Value = "M2"
$ws.Range("AA5").Value = "M1"
$ws.Range("AB5").Value = "M2"
$ws.Range("B6").Value = "A2"
$ws.Range("C6").Value = "DO"
$ws.Range("D6").Value = "M2"
$ws.Range("E6").Value = "A2"
$ws.Range("F6").Value = "M1"
$ws.Range("I6").Value = "M1"
$ws.Range("J6").Value = "DO"
$ws.Range("K6").Value = "M1"
$ws.Range("L6").Value = "M2"
$ws.Range("N6").Value = "M2"
$ws.Range("P6").Value = "M2"
$ws.Range("Q6").Value = "A1"
$ws.Range("R6").Value = "A1"
$ws.Range("T6").Value = "M2"
$ws.Range("V6").Value = "DO"
$ws.Range("W6").Value = "DO"
$ws.Range("X6").Value = "M2"
$ws.Range("Z6").Value = "A1"
$ws.Range("AA6").Value = "M2"
$ws.Range("AC6").Value = "A2"
$ws.Range("B7").Value = "DO"
$ws.Range("D7").Value = "A1"
$ws.Range("E7").Value = "M1"
$ws.Range("G7").Value = "A1"
$ws.Range("K7").Value = "DO"
$ws.Range("L7").Value = "A1"
$ws.Range("O7").Value = "M3"
$ws.Range("Q7").Value = "M3"
$ws.Range("U7").Value = "A1"
$ws.Range("V7").Value = "A1"
$ws.Range("Y7").Value = "M1"
$ws.Range("Z7").Value = "DO"
$ws.Range("AA7").Value = "A1"
$ws.Range("AC7").Value = "M3"
$ws.Range("B8").Value = "A2"
$ws.Range("H8").Value = "A1"
$ws.Range("I8").Value = "A1"
$ws.Range("J8").Value = "A2"
$ws.Range("K8").Value = "DO"
$ws.Range("L8").Value = "A2"
$ws.Range("N8").Value = "M2"
$ws.Range("O8").Value = "A2"
$ws.Range("P8").Value = "A2"
$ws.Range("R8").Value = "DO"
$ws.Range("T8").Value = "A2"
$ws.Range("V8").Value = "A2"
$ws.Range("W8").Value = "A2"
$ws.Range("Y8").Value = "DO"
$ws.Range("AA8").Value = "A1"
$ws.Range("AC8").Value = "A2"
$ws.Range("B9").Value = "M2"
$ws.Range("C9").Value = "A1"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = "DO"
$ws.Range("G9").Value = "A2"
$ws.Range("H9").Value = "M2"
$ws.Range("I9").Value = "M2"
$ws.Range("J9").Value = "M1"
$ws.Range("K9").Value = "M1"
$ws.Range("L9").Value = "M2"
$ws.Range("M9").Value = "M2"
$ws.Range("O9").Value = "DO"
$ws.Range("P9").Value = "M2"
$ws.Range("Q9").Value = "M2"
$ws.Range("T9").Value = "M1"
$ws.Range("U9").Value = "A2"
$ws.Range("V9").Value = "M1"
$ws.Range("W9").Value = "M2"
$ws.Range("X9").Value = "M2"
$ws.Range("Y9").Value = "M2"
$ws.Range("Z9").Value = "M1"
$ws.Range("AA9").Value = "M2"
$ws.Range("AB9").Value = "A1"
$ws.Range("AC9").Value = "DO"
$ws.Range("C10").Value = "M2"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = "M1"
$ws.Range("F10").Value = "A1"
$ws.Range("G10").Value = "M2"
$ws.Range("I10").Value = "M2"
$ws.Range("K10").Value = "M2"
$ws.Range("M10").Value = "A1"
$ws.Range("N10").Value = "M1"
$ws.Range("O10").Value = "DO"
$ws.Range("P10").Value = "M2"
$ws.Range("Q10").Value = "A2"
$ws.Range("R10").Value = "A1"
$ws.Range("T10").Value = "M2"
$ws.Range("U10").Value = "M1"
$ws.Range("W10").Value = "A2"
$ws.Range("X10").Value = "M2"
$ws.Range("Y10").Value = "M2"
$ws.Range("AA10").Value = "A1"
$ws.Range("AB10").Value = "M1"
